$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.275959610939026
$ws.Range("B1").Value = 1.923154950141907
$ws.Range("C1").Value = 5.59109354019165
$ws.Range("D1").Value = 1.933579325675964
$ws.Range("E1").Value = 1.113146185874939
